$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume/1h) store plain text (prices with
# dot-grouped digits, percentages with padding spaces), not numbers.
# Force each touched cell to a text number format *before* writing its
# new value so that numeric-looking strings (e.g. "1.02") stay text
# instead of being auto-converted to numbers, matching the original
# inline-string cell type. Only the cells that actually change are
# touched, so untouched neighboring cells keep their original format.
$changedCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "E26", "D27", "E27", "D28", "E28", "D29", "E29", "D30", "E30", "D31", "E31", "D32", "E32", "D33", "E33", "D34", "E34", "D35", "E35", "D36", "E36", "D37", "E37", "D38", "E38", "E39", "D40", "E40", "E41", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "E50", "D51", "E51")
foreach ($addr in $changedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "41.453.81"
$ws.Range("E2").Value = "  -1.22%  "
$ws.Range("D3").Value = "2.424.46"
$ws.Range("E3").Value = "  -2.67%  "
$ws.Range("D4").Value = "1.02"
$ws.Range("E4").Value = "  +1.55%  "
$ws.Range("D5").Value = "307.50"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").Value = "87.80"
$ws.Range("E6").Value = "  -8.45%  "
$ws.Range("D7").Value = "0.526"
$ws.Range("E7").Value = "  -5.69%  "
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("D9").Value = "0.474"
$ws.Range("E9").Value = "  -7.73%  "
$ws.Range("D10").Value = "31.43"
$ws.Range("E10").Value = "  -9.02%  "
$ws.Range("D11").Value = "0.0754"
$ws.Range("E11").Value = "  -4.56%  "
$ws.Range("D12").Value = "0.108"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "2.806.97"
$ws.Range("E13").Value = "  -2.52%  "
$ws.Range("D14").Value = "6.61"
$ws.Range("E14").Value = "  -6.73%  "
$ws.Range("D15").Value = "2.442.28"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "14.72"
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "0.749"
$ws.Range("E17").Value = "  -5.82%  "
$ws.Range("D18").Value = "41.072.53"
$ws.Range("E18").Value = "  -2.53%  "
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  -5.41%  "
$ws.Range("D20").Value = "0.0₃0885"
$ws.Range("E20").Value = "  -4.39%  "
$ws.Range("D21").Value = "68.14"
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("D22").Value = "10.56"
$ws.Range("E22").Value = "  -10.62%  "
$ws.Range("D23").Value = "229.24"
$ws.Range("E23").Value = "  -4.04%  "
$ws.Range("D24").Value = "2.63"
$ws.Range("E24").Value = "  -7.06%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -6.83%  "
$ws.Range("D27").Value = "23.20"
$ws.Range("E27").Value = "  -7.10%  "
$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").Value = "9.34"
$ws.Range("E29").Value = "  -5.11%  "
$ws.Range("D30").Value = "33.97"
$ws.Range("E30").Value = "  -8.44%  "
$ws.Range("D31").Value = "150.48"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").Value = "5.19"
$ws.Range("E32").Value = "  -9.40%  "
$ws.Range("D33").Value = "2.53"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("D34").Value = "2.49"
$ws.Range("E34").Value = "  -6.11%  "
$ws.Range("D35").Value = "0.0727"
$ws.Range("E35").Value = "  -5.44%  "
$ws.Range("D36").Value = "16.68"
$ws.Range("E36").Value = "  -3.47%  "
$ws.Range("D37").Value = "2.85"
$ws.Range("E37").Value = "  -6.39%  "
$ws.Range("D38").Value = "1.74"
$ws.Range("E38").Value = "  -8.84%  "
$ws.Range("E39").Value = "  -4.70%  "
$ws.Range("D40").Value = "0.0961"
$ws.Range("E40").Value = "  -8.74%  "
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("E42").Value = "  -5.78%  "
$ws.Range("D43").Value = "20.39"
$ws.Range("E43").Value = "  -4.47%  "
$ws.Range("D44").Value = "1.908.91"
$ws.Range("E44").Value = "  -5.04%  "
$ws.Range("D45").Value = "0.0270"
$ws.Range("E45").Value = "  -6.65%  "
$ws.Range("D46").Value = "2.81"
$ws.Range("E46").Value = "  -10.35%  "
$ws.Range("D47").Value = "2.686.42"
$ws.Range("E47").Value = "  -2.22%  "
$ws.Range("D48").Value = "8.41"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("D49").Value = "92.68"
$ws.Range("E49").Value = "  -5.84%  "
$ws.Range("E50").Value = "  -8.13%  "
$ws.Range("D51").Value = "70.59"
$ws.Range("E51").Value = "  -9.70%  "
